$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibition) - row => new F value
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 222
$ws1.Range("F4").Value = 12
$ws1.Range("F5").Value = 6630
$ws1.Range("F6").Value = 81
$ws1.Range("F8").Value = 134
$ws1.Range("F9").Value = 6052
$ws1.Range("F12").Value = 1240
$ws1.Range("F14").Value = 89
$ws1.Range("F16").Value = 111
$ws1.Range("F18").Value = 352
$ws1.Range("F19").Value = 42
$ws1.Range("F21").Value = 4358
$ws1.Range("F22").Value = 47
$ws1.Range("F23").Value = 16
$ws1.Range("F25").Value = 30

# Sheet "全部类型" (all types) - row => new F value
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 222
$ws4.Range("F4").Value = 12
$ws4.Range("F5").Value = 6630
$ws4.Range("F6").Value = 81
$ws4.Range("F8").Value = 134
$ws4.Range("F9").Value = 6052
$ws4.Range("F12").Value = 1240
$ws4.Range("F14").Value = 89
$ws4.Range("F16").Value = 111
$ws4.Range("F18").Value = 352
$ws4.Range("F19").Value = 42
$ws4.Range("F21").Value = 4358
$ws4.Range("F23").Value = 47
$ws4.Range("F24").Value = 16
$ws4.Range("F26").Value = 30
